# Apply cell-value updates from the cryptos list refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.028.01'
$ws.Range('E2').Value = '  -2.05%  '
$ws.Range('D3').Value = '2.975.29'
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.21%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  -1.38%  '
$ws.Range('D9').Value = '2.972.52'
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.98'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.73%  '
$ws.Range('E12').Value = '  +2.42%  '
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.02'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('D16').Value = '3.464.78'
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').Value = '61.149.70'
$ws.Range('E17').Value = '  -1.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.83'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.87%  '
$ws.Range('D19').Value = '2.972.91'
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '448.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.97'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.679'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.27'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.38'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.89'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.37%  '
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.67'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.04'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.23'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.108'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').Value = '0.0₃0805'
$ws.Range('E35').Value = '  +2.99%  '
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.76'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.96'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('E39').Value = '  -2.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.93'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('E41').Value = '  +6.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.82'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '386.98'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.267'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E45').Value = '  -0.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.47'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('D47').Value = '2.692.94'
$ws.Range('E47').Value = '  -2.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.59'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('E51').Value = '  -0.29%  '
